# Weekly update: insert a new price record for "Alcachofa" (Madrigal
# variety, Región del Maule) as row 103 in the "Vega Modelo de Temuco"
# sheet. All existing records from row 103 downward shift down by one
# row (103->104, 104->105, ..., 154->155); their data is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 103 down one row, leaving a blank row 103
# (this mirrors the cells below it, including number formats).
$ws.Rows(103).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A103").Value = 10
$ws.Range("B103").Value = "Vega Modelo de Temuco"
$ws.Range("C103").Value = "La Araucanía"
$ws.Range("D103").Value = 44518
$ws.Range("E103").Value = 9
$ws.Range("F103").Value = 100112013
$ws.Range("G103").Value = "Alcachofa"
$ws.Range("H103").Value = "Madrigal"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 110
$ws.Range("K103").Value = 12000
$ws.Range("L103").Value = 12000
$ws.Range("M103").Value = 12000
$ws.Range("N103").Value = "`$/caja 40 unidades"
$ws.Range("O103").Value = "Región del Maule"
$ws.Range("P103").Value = 300
$ws.Range("Q103").Value = 40
$ws.Range("R103").Value = "Hortaliza"
